# Slide 1 has two multi-run paragraphs (one run per word/space, as PowerPoint
# sometimes produces) that should be collapsed into a single run each, with
# the same overall text. Setting TextRange.Text to the already-matching
# concatenation is a no-op for some engines, so first nudge the value to
# something different ("x") to force the paragraph to be rebuilt as one run,
# then assign the real target text.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "x"
$title.Text = "A Table, with a caption"

$caption = $s.Shapes.Item(3).TextFrame.TextRange
$caption.Text = "x"
$caption.Text = "Demonstration of simple table syntax, with alignment"
